# Add a new "10-10-2020" column (column Y) to the COVID19 time-series sheet,
# mirroring the existing layout used by the preceding date columns (N:X).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell Y1: "10-10-2020" (stored as text, bold + bordered like X1) ---
$header = $ws.Cells.Item(1, 25)
$header.Value = "'10-10-2020"
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous (thin box border)

# --- Data rows 2-36: new daily totals for column Y ---
$values = @{
    2  = 3724
    3  = 691040
    4  = 8877
    5  = 161904
    6  = 182121
    7  = 11505
    8  = 108935
    9  = 3037
    10 = 276046
    11 = 32317
    12 = 129304
    13 = 127540
    14 = 13876
    15 = 69979
    16 = 81654
    17 = 561610
    18 = 175304
    19 = 3886
    20 = 124887
    21 = 1229339
    22 = 9866
    23 = 4903
    24 = 1974
    25 = 5656
    26 = 220388
    27 = 25543
    28 = 108533
    29 = 131766
    30 = 2721
    31 = 591811
    32 = 183025
    33 = 23801
    34 = 46058
    35 = 383086
    36 = 252806
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 25).Value = $values[$row]
}
